$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all bio for ELC")

$ws.Range("C13").Value = "~TFM_INS"

$ws.Range("C14").Value = "attribute"
$ws.Range("D14").Value = "pset_pn"
$ws.Range("E14").Value = "cset_cn"
$ws.Range("F14").Value = "value"

$ws.Range("C15").Value = "flo_cost"
$ws.Range("D15").Value = "IMPNRGZ"
$ws.Range("E15").Value = "UC_all_bio*"
$ws.Range("F15").Value = -100

$ws.Range("D17").Select()
